$d = $word.ActiveDocument

# The ordered list of old/new values taken from the diff, in document order.
$replacements = @(
    @{ Old = "60÷2="; New = "45÷5=" },
    @{ Old = "64÷9="; New = "48÷7=" },
    @{ Old = "72÷8="; New = "86÷8=" },
    @{ Old = "72÷8="; New = "38÷9=" },
    @{ Old = "57÷6="; New = "62÷9=" },
    @{ Old = "30÷4="; New = "71÷2=" },
    @{ Old = "97÷6="; New = "83÷7=" },
    @{ Old = "72÷7="; New = "11÷2=" },
    @{ Old = "34÷9="; New = "35÷4=" },
    @{ Old = "91÷2="; New = "92÷9=" },
    @{ Old = "67÷9="; New = "10÷8=" },
    @{ Old = "49÷2="; New = "72÷2=" },
    @{ Old = "53÷7="; New = "32÷3=" },
    @{ Old = "92÷2="; New = "99÷7=" },
    @{ Old = "59÷7="; New = "67÷5=" },
    @{ Old = "35÷2="; New = "73÷3=" },
    @{ Old = "78÷2="; New = "68÷6=" },
    @{ Old = "34÷5="; New = "57÷2=" },
    @{ Old = "21÷2="; New = "38÷5=" },
    @{ Old = "84÷7="; New = "42÷4=" },
    @{ Old = "87÷5="; New = "19÷9=" },
    @{ Old = "30÷8="; New = "37÷5=" },
    @{ Old = "47÷9="; New = "59÷7=" },
    @{ Old = "19÷5="; New = "15÷8=" },
    @{ Old = "59÷4="; New = "28÷8=" }
)

# Find every paragraph that holds one of these "a÷b=" expressions and
# remember its character span. Using Paragraphs()/Cell() objects directly
# to drive Find.Execute can mis-target when two paragraphs share identical
# text, so instead we resolve stable (Start, End) character offsets up
# front and replace through $d.Range(start, end), which addresses the
# story by position rather than by (possibly ambiguous) content.
$spans = New-Object System.Collections.ArrayList
$total = $d.Paragraphs.Count
for ($i = 1; $i -le $total; $i++) {
    $p = $d.Paragraphs($i)
    $r = $p.Range
    if ($r.Text -match "÷") {
        [void]$spans.Add(@{ Start = $r.Start; End = $r.End })
    }
}

if ($spans.Count -ne $replacements.Count) {
    throw "Expected $($replacements.Count) math expressions, found $($spans.Count)"
}

for ($idx = 0; $idx -lt $replacements.Count; $idx++) {
    $span = $spans[$idx]
    $rep = $replacements[$idx]
    $target = $d.Range($span.Start, $span.End)
    $found = $target.Find.Execute($rep.Old, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $rep.New, 2)
    if (-not $found) {
        throw "Could not find '$($rep.Old)' at span $($span.Start)-$($span.End)"
    }
}

Write-Host "Applied $($replacements.Count) cell updates."
